# Applies the "Automatic update of files" commit:
#  1. Bumps the "Förändrad" date (column C) from 46062 -> 46063 for every
#     data row (rows 2-24).
#  2. Re-orders the data rows 10-24 (same 15 records, shuffled into a new
#     row order) while every other row (2-9) keeps its own data in place.
#
# Because the destination/source ranges for the re-order overlap, we take a
# full snapshot (values/formulas) of rows 10-24 with COM *before* writing
# anything back, then re-emit each destination row from its mapped source
# row. To keep the sheet as close as possible to its original state (and
# avoid needless float round-tripping through COM) we only touch a cell
# when its target content actually differs from what is already sitting
# in that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that are ever populated in the data rows.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

function Get-CellSnapshot($cell) {
    if ($cell.HasFormula) {
        return @{ kind = "formula"; value = $cell.Formula }
    }
    $v = $cell.Value2
    if ($v -eq $null) {
        return @{ kind = "empty"; value = $null }
    }
    return @{ kind = "value"; value = $v }
}

# ---------------------------------------------------------------------
# Step 1: bump column C (Förändrad) for every data row, rows 2-24.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 24; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cur = $cCell.Value2
    if ($cur -eq 46062) {
        $cCell.Value2 = 46063
    }
}

# ---------------------------------------------------------------------
# Step 2: snapshot rows 10-24 (post column-C bump) before reshuffling.
# ---------------------------------------------------------------------
$snapshot = @{}
for ($r = 10; $r -le 24; $r++) {
    $rowData = @{}
    foreach ($col in $columns) {
        $rowData[$col] = Get-CellSnapshot $ws.Range($col + $r)
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (within the pre-reorder snapshot above).
$rowMap = @{10=14; 11=10; 12=13; 13=11; 14=12; 15=15; 16=16; 17=21; 18=17; 19=18; 20=22; 21=24; 22=23; 23=20; 24=19}

# ---------------------------------------------------------------------
# Step 3: write each destination row from its mapped source snapshot,
# but only where the content actually changes, and clear any column the
# destination no longer needs.
# ---------------------------------------------------------------------
for ($destRow = 10; $destRow -le 24; $destRow++) {
    $srcRow = $rowMap[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $srcData = $snapshot[$srcRow]
    foreach ($col in $columns) {
        $cell = $ws.Range($col + $destRow)
        $info = $srcData[$col]
        $current = Get-CellSnapshot $cell

        if ($info.kind -eq "formula") {
            if ($current.kind -ne "formula" -or $current.value -ne $info.value) {
                $cell.Formula = $info.value
            }
        } elseif ($info.kind -eq "value") {
            if ($current.kind -ne "value" -or $current.value -ne $info.value) {
                $cell.Value2 = $info.value
            }
        } else {
            if ($current.kind -ne "empty") {
                $cell.ClearContents()
            }
        }
    }
}
